# 10.1.1 sheet refresh: replace the 2007-2019 series (columns D:P) with the
# new 2015-2021 series (columns D:J), dropping the now-unused trailing years.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix up the formatting of the cells that will carry the new data ---
# D5:G5 currently use the "old" number-format style (the one that the diff
# removes from cellXfs); re-point them at the style already used by H5 so
# that style becomes unused everywhere in the sheet.
$ws.Range("H5").Copy()
$ws.Range("D5:G5").PasteSpecial(-4122)

# E6:L6 currently share D6's style; the new layout wants them to match the
# style used further right (M6), with its thinner top border.
$ws.Range("M6").Copy()
$ws.Range("E6:L6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 2. Write the new header years (row 4) ---
$ws.Range("D4").Value = 2015
$ws.Range("E4").Value = 2016
$ws.Range("F4").Value = 2017
$ws.Range("G4").Value = 2018
$ws.Range("H4").Value = 2019
$ws.Range("I4").Value = 2020
$ws.Range("J4").Value = 2021

# --- 3. Write the new "bottom 40%" series (row 5) ---
$ws.Range("D5").Value = 2.2197193775563164
$ws.Range("E5").Value = 2.1235271668715399
$ws.Range("F5").Value = 2.7818537161298167
$ws.Range("G5").Value = 6.7272960584548969
$ws.Range("H5").Value = 5.1525830614767187
$ws.Range("I5").Value = 4.4774536255935971
$ws.Range("J5").Value = 4.6024666695867751

# --- 4. Write the new "total population" series (row 6) ---
$ws.Range("D6").Value = 2.2322863217945752
$ws.Range("E6").Value = 2.8603553109638966
$ws.Range("F6").Value = 3.113207036164539
$ws.Range("G6").Value = 6.2970593463100784
$ws.Range("H6").Value = 4.8617746111834492
$ws.Range("I6").Value = 2.6715092780025032
$ws.Range("J6").Value = 4.3694509108608912

# --- 5. Wipe the now-obsolete trailing year columns (K:P) so the sheet's
# used range shrinks back to A1:J6, without reshuffling the 16384-column grid
# the way a real column delete would. ---
$ws.Range("K1:P6").Clear()

# --- 6. Match the new column widths for D:J and the recorded selection ---
$ws.Range("D1:J1").EntireColumn.ColumnWidth = 8.67
$ws.Range("K16").Select()
